$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column G with header and data for 07_02_2024
$ws.Range("G1").Value = "07_02_2024"
$ws.Range("G2").Value = 1097
$ws.Range("G3").Value = 1072
$ws.Range("G4").Value = 1157
$ws.Range("G5").Value = 2930

# Update selection to match the diff (active cell G5)
[void]$ws.Range("G5").Select()
